$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: G. d. Arrascaeta / Flamengo -> L. Mugni / Ceara (assists stay 5)
$ws.Range("A2").Value = "L. Mugni"
$ws.Range("B2").Value = "Ceara"
$ws.Range("C2").Value = 5

# Row 3: M. Pereira / Cruzeiro / 5 (unchanged)
$ws.Range("A3").Value = "M. Pereira"
$ws.Range("B3").Value = "Cruzeiro"
$ws.Range("C3").Value = 5

# Row 4: A. Patrick / Internacional -> K. Jorge / Cruzeiro (assists stay 5)
$ws.Range("A4").Value = "K. Jorge"
$ws.Range("B4").Value = "Cruzeiro"
$ws.Range("C4").Value = 5

# Row 5: K. Jorge / Cruzeiro / 4 -> A. Patrick / Internacional / 5
$ws.Range("A5").Value = "A. Patrick"
$ws.Range("B5").Value = "Internacional"
$ws.Range("C5").Value = 5

# Row 6: J. Arias / Fluminense / 4 -> G. d. Arrascaeta / Flamengo / 5
$ws.Range("A6").Value = "G. d. Arrascaeta"
$ws.Range("B6").Value = "Flamengo"
$ws.Range("C6").Value = 5
